$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135, shifting existing rows 135-182 down to 136-183.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new weekly data point.
$ws.Cells.Item(135, 1).Value = 4
$ws.Cells.Item(135, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(135, 3).Value = "Los Lagos"
$ws.Cells.Item(135, 4).Value = 44463
$ws.Cells.Item(135, 5).Value = 10
$ws.Cells.Item(135, 6).Value = 100114014
$ws.Cells.Item(135, 7).Value = "Betarraga"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 1200
$ws.Cells.Item(135, 11).Value = 1200
$ws.Cells.Item(135, 12).Value = 1200
$ws.Cells.Item(135, 13).Value = 1200
$ws.Cells.Item(135, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(135, 15).Value = "Región del Maule"
$ws.Cells.Item(135, 16).Value = 240
$ws.Cells.Item(135, 17).Value = 5
$ws.Cells.Item(135, 18).Value = "Hortaliza"
